$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) — sheet1
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 10211
$ws1.Range("F9").Value = 771
$ws1.Range("F11").Value = 1222
$ws1.Range("F12").Value = 1061
$ws1.Range("F13").Value = 3179
$ws1.Range("F14").Value = 2379
$ws1.Range("F16").Value = 2106
$ws1.Range("F17").Value = 2106
$ws1.Range("F18").Value = 245
$ws1.Range("F34").Value = 243
$ws1.Range("F37").Value = 402
$ws1.Range("F38").Value = 395
$ws1.Range("F39").Value = 1684
$ws1.Range("F41").Value = 429
$ws1.Range("F44").Value = 975
$ws1.Range("F46").Value = 357

# Sheet "演出" (Performance) — sheet2
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 40

# Sheet "全部类型" (All Types) — sheet4
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 10211
$ws4.Range("F11").Value = 771
$ws4.Range("F12").Value = 1061
$ws4.Range("F13").Value = 3179
$ws4.Range("F14").Value = 2379
$ws4.Range("F15").Value = 2106
$ws4.Range("F16").Value = 2106
$ws4.Range("F17").Value = 245
$ws4.Range("F32").Value = 40
$ws4.Range("F36").Value = 243
$ws4.Range("F39").Value = 404
$ws4.Range("F41").Value = 395
$ws4.Range("F42").Value = 1684
$ws4.Range("F45").Value = 429
$ws4.Range("F48").Value = 975
$ws4.Range("F49").Value = 357
